$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Support Vector Machines (SVM)*") {
        $p.Range.Font.StrikeThrough = 1
    }
    if ($t -like "LightGBM: LightGBM is a gradient boosting framework*") {
        $p.Range.Font.StrikeThrough = 1
    }
}
